$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 8 new rows before row 19, shifting existing rows 19-121 down to 27-129
$ws.Rows("19:26").Insert()

# Format the new data cells (B:F) as Text so numeric-looking strings keep their exact formatting
$ws.Range("B19:F26").NumberFormat = "@"

# Row 19: Common Plus Preferred Stock
$ws.Range("A19").Value = 'Common Plus Preferred Stock'
$ws.Range("B19").Value = '338,752,830.58'
$ws.Range("C19").Value = '-74,100.00'
$ws.Range("D19").Value = '6,000,000.00'
$ws.Range("E19").Value = '9,817,134,000.00'
$ws.Range("F19").Value = '926,606,262.05'
$ws.Range("G19").Value = 'Constructed for Altman''s Z'

# Row 20: EBIT
$ws.Range("A20").Value = 'EBIT'
$ws.Range("B20").Value = '304,424,183.99'
$ws.Range("C20").Value = '-1,364,004,000.00'
$ws.Range("D20").Value = '122,944,000.00'
$ws.Range("E20").Value = '4,334,000,000.00'
$ws.Range("F20").Value = '556,409,802.60'
$ws.Range("G20").Value = 'Constructed for Altman''s Z'

# Row 21: Ratio A
$ws.Range("A21").Value = 'Ratio A'
$ws.Range("B21").Value = '0.02'
$ws.Range("C21").Value = '-0.02'
$ws.Range("D21").Value = '0.02'
$ws.Range("E21").Value = '0.08'
$ws.Range("F21").Value = '0.02'
$ws.Range("G21").Value = 'Constructed for Altman''s Z'

# Row 22: Ratio B
$ws.Range("A22").Value = 'Ratio B'
$ws.Range("B22").Value = '0.21'
$ws.Range("C22").Value = '0.04'
$ws.Range("D22").Value = '0.18'
$ws.Range("E22").Value = '0.70'
$ws.Range("F22").Value = '0.15'
$ws.Range("G22").Value = 'Constructed for Altman''s Z'

# Row 23: Ratio C
$ws.Range("A23").Value = 'Ratio C'
$ws.Range("B23").Value = '1.87'
$ws.Range("C23").Value = '0.29'
$ws.Range("D23").Value = '1.42'
$ws.Range("E23").Value = '8.06'
$ws.Range("F23").Value = '1.55'
$ws.Range("G23").Value = 'Constructed for Altman''s Z'

# Row 24: Ratio D
$ws.Range("A24").Value = 'Ratio D'
$ws.Range("B24").Value = '0.13'
$ws.Range("C24").Value = '-0.13'
$ws.Range("D24").Value = '0.10'
$ws.Range("E24").Value = '0.57'
$ws.Range("F24").Value = '0.15'
$ws.Range("G24").Value = 'Constructed for Altman''s Z'

# Row 25: Ratio E
$ws.Range("A25").Value = 'Ratio E'
$ws.Range("B25").Value = '0.23'
$ws.Range("C25").Value = '-0.74'
$ws.Range("D25").Value = '0.21'
$ws.Range("E25").Value = '1.01'
$ws.Range("F25").Value = '0.32'
$ws.Range("G25").Value = 'Constructed for Altman''s Z'

# Row 26: Working Capital
$ws.Range("A26").Value = 'Working Capital'
$ws.Range("B26").Value = '1,125,108,587.77'
$ws.Range("C26").Value = '-28,931,855,000.00'
$ws.Range("D26").Value = '543,614,000.00'
$ws.Range("E26").Value = '39,464,552,600.00'
$ws.Range("F26").Value = '3,845,915,891.90'
$ws.Range("G26").Value = 'Constructed for Altman''s Z'
